$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Greg")
$ws.Range("H13").Value = ""
$ws.Range("F13").Value = 8
$ws.Range("H14").Value = 5
